{"js": "// Replace two-digit division problems per the diff.\nconst replacements = [\n  [\"71\u00f72=35, 1\", \"63\u00f78=7, 7\"],\n  [\"74\u00f74=18, 2\", \"53\u00f77=7, 4\"],\n  [\"54\u00f79=6, 0\", \"14\u00f78=1, 6\"],\n  [\"93\u00f75=18, 3\", \"45\u00f76=7, 3\"],\n  [\"24\u00f76=4, 0\", \"91\u00f73=30, 1\"],\n  [\"80\u00f77=11, 3\", \"17\u00f77=2, 3\"],\n  [\"76\u00f76=12, 4\", \"34\u00f72=17, 0\"],\n  [\"25\u00f76=4, 1\", \"77\u00f79=8, 5\"],\n  [\"98\u00f73=32, 2\", \"82\u00f74=20, 2\"],\n  [\"17\u00f74=4, 1\", \"49\u00f78=6, 1\"],\n  [\"50\u00f73=16, 2\", \"62\u00f72=31, 0\"],\n  [\"90\u00f75=18, 0\", \"75\u00f74=18, 3\"],\n  [\"15\u00f72=7, 1\", \"65\u00f76=10, 5\"],\n  [\"44\u00f78=5, 4\", \"76\u00f79=8, 4\"],\n  [\"20\u00f74=5, 0\", \"91\u00f73=30, 1\"],\n  [\"39\u00f75=7, 4\", \"20\u00f75=4, 0\"],\n  [\"36\u00f73=12, 0\", \"74\u00f77=10, 4\"],\n  [\"35\u00f72=17, 1\", \"58\u00f73=19, 1\"],\n  [\"37\u00f74=9, 1\", \"87\u00f79=9, 6\"],\n  [\"88\u00f78=11, 0\", \"20\u00f78=2, 4\"],\n  [\"86\u00f72=43, 0\", \"18\u00f72=9, 0\"],\n  [\"29\u00f78=3, 5\", \"90\u00f76=15, 0\"],\n  [\"64\u00f72=32, 0\", \"24\u00f77=3, 3\"],\n  [\"23\u00f74=5, 3\", \"12\u00f73=4, 0\"],\n  [\"60\u00f74=15, 0\", \"67\u00f72=33, 1\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace two-digit division problems per the diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"71\u00f72=35, 1\", \"63\u00f78=7, 7\"),\n    @(\"74\u00f74=18, 2\", \"53\u00f77=7, 4\"),\n    @(\"54\u00f79=6, 0\", \"14\u00f78=1, 6\"),\n    @(\"93\u00f75=18, 3\", \"45\u00f76=7, 3\"),\n    @(\"24\u00f76=4, 0\", \"91\u00f73=30, 1\"),\n    @(\"80\u00f77=11, 3\", \"17\u00f77=2, 3\"),\n    @(\"76\u00f76=12, 4\", \"34\u00f72=17, 0\"),\n    @(\"25\u00f76=4, 1\", \"77\u00f79=8, 5\"),\n    @(\"98\u00f73=32, 2\", \"82\u00f74=20, 2\"),\n    @(\"17\u00f74=4, 1\", \"49\u00f78=6, 1\"),\n    @(\"50\u00f73=16, 2\", \"62\u00f72=31, 0\"),\n    @(\"90\u00f75=18, 0\", \"75\u00f74=18, 3\"),\n    @(\"15\u00f72=7, 1\", \"65\u00f76=10, 5\"),\n    @(\"44\u00f78=5, 4\", \"76\u00f79=8, 4\"),\n    @(\"20\u00f74=5, 0\", \"91\u00f73=30, 1\"),\n    @(\"39\u00f75=7, 4\", \"20\u00f75=4, 0\"),\n    @(\"36\u00f73=12, 0\", \"74\u00f77=10, 4\"),\n    @(\"35\u00f72=17, 1\", \"58\u00f73=19, 1\"),\n    @(\"37\u00f74=9, 1\", \"87\u00f79=9, 6\"),\n    @(\"88\u00f78=11, 0\", \"20\u00f78=2, 4\"),\n    @(\"86\u00f72=43, 0\", \"18\u00f72=9, 0\"),\n    @(\"29\u00f78=3, 5\", \"90\u00f76=15, 0\"),\n    @(\"64\u00f72=32, 0\", \"24\u00f77=3, 3\"),\n    @(\"23\u00f74=5, 3\", \"12\u00f73=4, 0\"),\n    @(\"60\u00f74=15, 0\", \"67\u00f72=33, 1\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $null = $find.Execute(\n        $oldText,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $newText,\n        2\n    )\n}\n"}
